$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '87.566.86'
$ws.Range("E2").Value = '  +3.24%  '

$ws.Range("D3").Value = '3.251.07'
$ws.Range("E3").Value = '  -1.99%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.79%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '626.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.72%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.378'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +16.93%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.687'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +16.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("D10").Value = '3.250.78'
$ws.Range("E10").Value = '  -2.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.576'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.183'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +9.96%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000259'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.92%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.861.57'
$ws.Range("E14").Value = '  -1.45%  '

$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '33.91'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.32'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.01%  '

$ws.Range("D17").Value = '87.336.87'
$ws.Range("E17").Value = '  +3.16%  '

$ws.Range("D18").Value = '3.252.76'
$ws.Range("E18").Value = '  -1.43%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.71%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.99'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.69%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '432.71'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.39%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.86'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.98%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.33'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.82%  '

$ws.Range("E24").Value = '  -1.50%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.40'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.79%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.11'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.61%  '

$ws.Range("D27").Value = '3.411.37'
$ws.Range("E27").Value = '  -2.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '76.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.14%  '

$ws.Range("E29").Value = '  -0.60%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.14%  '

$ws.Range("E31").Value = '  +10.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.29%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.77'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '551.73'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.57%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.39'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -12.62%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.95'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.97'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.45%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.137'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -10.40%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.42'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.82%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '21.72'
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.392'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.45%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.93'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.46%  '

$ws.Range("E45").Value = '  -0.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '152.51'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.27%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '178.87'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '44.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.30%  '

$ws.Range("E49").Value = '  -4.68%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.22'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.65%  '

$ws.Range("E51").Value = '  +10.70%  '

